$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DSLS"

# Row 1 (header) - reuse existing styled cells A1:E1
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Quyen"
$ws.Range("D1").Value = "Matkhau"
$ws.Range("E1").Value = "Ngaytao"

# Row 2 - reuse existing styled cells A2:D2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "i@gmail.com"
$ws.Range("C2").Value = "Nhân viên"
$ws.Range("D2").Value = "anhkk"

# Row 3 - reuse existing styled cells A3:D3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "email"
$ws.Range("C3").Value = "Khách hàng"
$ws.Range("D3").Value = "mk"

# Row 4 - new row, copy style from row 3 first
$ws.Range("A3:D3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = 993
$ws.Range("B4").Value = "quynh181204@gmail.com"
$ws.Range("C4").Value = "Khách hàng"
$ws.Range("D4").Value = "quynh"

# Remove old leftover cells in E2,E3 and F:G columns
$ws.Range("E2:G3").Clear()
$ws.Range("F1:G1").Clear()

# E1 becomes style "Normal" (s=0) instead of header style
$ws.Range("Z1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
